$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the abstract (column D, row 6) with the cleaned-up text for the
# "Nasopharyngeal Microbiota Profiling of SARS-CoV-2 Infected Patients" row.
$abstract = "We analyzed the bacterial communities of the nasopharynx in 40 SARS-CoV-2 infected and uninfected patients.`n All infected patients had a mild COVID-19 disease.`n We did not find statistically significant differences in either bacterial richness and diversity or composition.`n These findings suggest a nasopharyngeal microbiota at least early resilient to SARS-CoV-2 infection.`n"
$ws.Range("D6").Value = $abstract

# Update the authors list (column E, row 6) - spacing between authors
# normalized to four spaces instead of three.
$authors = "[Flavio%De Maio%NULL%1,    Brunella%Posteraro%NULL%1,    Francesca Romana%Ponziani%NULL%1,    Paola%Cattani%NULL%1,    Antonio%Gasbarrini%NULL%0,    Maurizio%Sanguinetti%maurizio.sanguinetti@unicatt.it%1]"
$ws.Range("E6").Value = $authors
